$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new feedback for Breast Cancer
$ws.Range("A5").Value = 697696
$ws.Range("B5").Value = "Patient"
$ws.Range("D5").Value = "X"
$ws.Range("E5").Value2 = 45638
$ws.Range("F5").Value = "Breast Cancer"
$ws.Range("C5").Value = "Should talk more about male breast though"

# Row 6: new feedback for Kidney Cancer
$ws.Range("A6").Value = 697696
$ws.Range("B6").Value = "Patient"
$ws.Range("D6").Value = "X"
$ws.Range("E6").Value2 = 45330
$ws.Range("F6").Value = "Kidney Cancer"
$ws.Range("C6").Value = "WTF is this"

# Match the date number format used by the existing date column (E2:E4)
$ws.Range("E2").Copy()
$ws.Range("E5:E6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C6").Select() | Out-Null
